$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header from total_shared_rides to avg_shared_rides
$ws.Range("C1").Value = "avg_shared_rides"

# Update values in column C to reflect averages instead of totals
$ws.Range("C2").Value = 5046
$ws.Range("C3").Value = 4724.1
$ws.Range("C4").Value = -6.379310344827579
$ws.Range("C5").Value = 3763.396226415094
$ws.Range("C6").Value = 3720.830188679245
$ws.Range("C7").Value = -1.131053845382532
